$d = $word.ActiveDocument

$rng = $d.Content
[void]$rng.Find.Execute("Research & Data Analytics Leadership", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$bullet = [char]0x2022
$newText = "`r" + $bullet + " Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters" `
    + "`r" + $bullet + " Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States" `
    + "`r" + $bullet + " Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"

[void]$rng.InsertAfter($newText)
